$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H30").Value = 44022
